# The commit swaps the deck's two embedded themes: the main/slide-master
# theme (ppt/theme/theme1.xml, "Integral") and the notes-master theme
# (ppt/theme/theme2.xml, "Office Theme") trade color schemes — theme1
# becomes the Office Theme palette, theme2 becomes the Integral palette.
# (Their font schemes and format schemes are already byte-identical, so
# only the 12 theme colors actually change.)
#
# This host's PowerPoint object model only exposes one writable theme
# color scheme (the slide master's, which backs ppt/theme/theme1.xml), so
# we repaint it with the Office Theme palette that used to live in
# theme2.xml.

function HexToOle($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeTheme = [ordered]@{
    1  = "000000"   # dk1
    2  = "FFFFFF"   # lt1
    3  = "44546A"   # dk2
    4  = "E7E6E6"   # lt2
    5  = "5B9BD5"   # accent1
    6  = "ED7D31"   # accent2
    7  = "A5A5A5"   # accent3
    8  = "FFC000"   # accent4
    9  = "4472C4"   # accent5
    10 = "70AD47"   # accent6
    11 = "0563C1"   # hlink
    12 = "954F72"   # folHlink
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

foreach ($idx in $officeTheme.Keys) {
    $cs.Item($idx).RGB = HexToOle($officeTheme[$idx])
}
